$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CAL")

# Row 21: Earnings Before Interest And Taxes -> J21 becomes "NA"
$ws.Range("J21").Value = "NA"

# Row 83: Depreciation -> J83 becomes "NA"
$ws.Range("J83").Value = "NA"

# Row 91: Capital Expenditures -> update D:J values
$ws.Range("D91").Value = -44700
$ws.Range("E91").Value = -50500
$ws.Range("F91").Value = -73500
$ws.Range("G91").Value = -45000
$ws.Range("H91").Value = -44000
$ws.Range("I91").Value = -55800
$ws.Range("J91").Value = -27900

# Row 94: Total Cash Flows From Investing Activities -> J94 becomes "NA"
$ws.Range("J94").Value = "NA"

# Row 100: Total Cash Flows From Financing Activities -> J100 becomes "NA"
$ws.Range("J100").Value = "NA"

# Row 101: Effect Of Exchange Rate Changes -> J101 becomes "NA"
$ws.Range("J101").Value = "NA"
